$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item("IND_Regression_FullCycle")
$ws4 = $wb.Worksheets.Item("CHN_Regression_FullCycle")

# ---------------------------------------------------------------------
# Remove the now-obsolete trailing rows (old rows 15:33) from both
# sheets first, while the two new "Block Global" rows still need to be
# written in place of rows 2:3 (the old rows 2:12 simply get overwritten
# with the data that used to live one row higher / two rows higher, see
# below) - this keeps every surviving cell's original style (s="2") and
# avoids the formatting bleed that Rows.Insert() would introduce.
# ---------------------------------------------------------------------
$ws3.Rows("15:33").Delete()
$ws4.Rows("15:33").Delete()

# ---------------------------------------------------------------------
# IND_Regression_FullCycle (sheet3): rows 2:14
# ---------------------------------------------------------------------
$ws3.Range("A2").Value = "BlockGlobalClient"
$ws3.Range("B2").Value = "Blockglobalclient"
$ws3.Range("C2").Value = "Block created global client"
$ws3.Range("D2").Value = "No"
$ws3.Range("A3").Value = "BlockGlobalBrand"
$ws3.Range("B3").Value = "Blockglobalbrand"
$ws3.Range("C3").Value = "Block created global brand"
$ws3.Range("D3").Value = "Yes"
$ws3.Range("A4").Value = "JobCreation"
$ws3.Range("B4").Value = "createJob"
$ws3.Range("C4").Value = "Create Job for Opco"
$ws3.Range("D4").Value = "No"
$ws3.Range("A5").Value = "CreateSubJob"
$ws3.Range("B5").Value = "SubJob"
$ws3.Range("C5").Value = "Create Sub Job for Created Job"
$ws3.Range("D5").Value = "Yes"
$ws3.Range("A6").Value = "CreateBudget"
$ws3.Range("B6").Value = "createBudget"
$ws3.Range("C6").Value = "Create Working Estimate for Job"
$ws3.Range("D6").Value = "No"
$ws3.Range("A7").Value = "CreateQuote"
$ws3.Range("B7").Value = "CreateQuote"
$ws3.Range("C7").Value = "Create Quote"
$ws3.Range("D7").Value = "No"
$ws3.Range("A8").Value = "FixedAssetPurchaseOrder"
$ws3.Range("B8").Value = "CreatePurchaseOrder"
$ws3.Range("C8").Value = "Create Purchase Order"
$ws3.Range("D8").Value = "No"
$ws3.Range("A9").Value = "CreatePurchaseOrder"
$ws3.Range("B9").Value = "CreatePurchaseOrder"
$ws3.Range("C9").Value = "Create Purchase Order"
$ws3.Range("D9").Value = "No"
$ws3.Range("A10").Value = "ApprovePurchaseOrder"
$ws3.Range("B10").Value = "ApprovePurchaseOrder"
$ws3.Range("C10").Value = "Approve Created PurchaseOrder"
$ws3.Range("D10").Value = "No"
$ws3.Range("A11").Value = "RejectPurchaseOrder"
$ws3.Range("B11").Value = "RejectPurchaseOrder"
$ws3.Range("C11").Value = "Reject Created PurchaseOrder"
$ws3.Range("D11").Value = "No"
$ws3.Range("A12").Value = "CreateVendorInvoice"
$ws3.Range("B12").Value = "CreateInvoice"
$ws3.Range("C12").Value = "Create invoice for purchase order"
$ws3.Range("D12").Value = "No"
$ws3.Range("A13").Value = "ApproveVendorInvoice"
$ws3.Range("B13").Value = "ApproveInvoice"
$ws3.Range("C13").Value = "Approve Created Vendor Invoice"
$ws3.Range("D13").Value = "No"
$ws3.Range("A14").Value = "RejectVendorInvoice"
$ws3.Range("B14").Value = "RejectInvoice"
$ws3.Range("C14").Value = "Reject Created Vendor Invoice"
$ws3.Range("D14").Value = "No"

# ---------------------------------------------------------------------
# CHN_Regression_FullCycle (sheet4): rows 2:14
# ---------------------------------------------------------------------
$ws4.Range("A2").Value = "BlockGlobalClient"
$ws4.Range("B2").Value = "Blockglobalclient"
$ws4.Range("C2").Value = "Block created global client"
$ws4.Range("D2").Value = "Yes"
$ws4.Range("A3").Value = "BlockGlobalBrand"
$ws4.Range("B3").Value = "Blockglobalbrand"
$ws4.Range("C3").Value = "Block created global brand"
$ws4.Range("D3").Value = "No"
$ws4.Range("A4").Value = "JobCreation"
$ws4.Range("B4").Value = "createJob"
$ws4.Range("C4").Value = "Create Job for Opco"
$ws4.Range("D4").Value = "No"
$ws4.Range("A5").Value = "CreateSubJob"
$ws4.Range("B5").Value = "SubJob"
$ws4.Range("C5").Value = "Create Sub Job for Created Job"
$ws4.Range("D5").Value = "Yes"
$ws4.Range("A6").Value = "CreateBudget"
$ws4.Range("B6").Value = "createBudget"
$ws4.Range("C6").Value = "Create Working Estimate for Job"
$ws4.Range("D6").Value = "No"
$ws4.Range("A7").Value = "CreateQuote"
$ws4.Range("B7").Value = "CreateQuote"
$ws4.Range("C7").Value = "Create Quote"
$ws4.Range("D7").Value = "No"
$ws4.Range("A8").Value = "FixedAssetPurchaseOrder"
$ws4.Range("B8").Value = "CreatePurchaseOrder"
$ws4.Range("C8").Value = "Create Purchase Order"
$ws4.Range("D8").Value = "No"
$ws4.Range("A9").Value = "CreatePurchaseOrder"
$ws4.Range("B9").Value = "CreatePurchaseOrder"
$ws4.Range("C9").Value = "Create Purchase Order"
$ws4.Range("D9").Value = "No"
$ws4.Range("A10").Value = "ApprovePurchaseOrder"
$ws4.Range("B10").Value = "ApprovePurchaseOrder"
$ws4.Range("C10").Value = "Approve Created PurchaseOrder"
$ws4.Range("D10").Value = "No"
$ws4.Range("A11").Value = "RejectPurchaseOrder"
$ws4.Range("B11").Value = "RejectPurchaseOrder"
$ws4.Range("C11").Value = "Reject Created PurchaseOrder"
$ws4.Range("D11").Value = "No"
$ws4.Range("A12").Value = "CreateVendorInvoice"
$ws4.Range("B12").Value = "CreateInvoice"
$ws4.Range("C12").Value = "Create invoice for purchase order"
$ws4.Range("D12").Value = "No"
$ws4.Range("A13").Value = "ApproveVendorInvoice"
$ws4.Range("B13").Value = "ApproveInvoice"
$ws4.Range("C13").Value = "Approve Created Vendor Invoice"
$ws4.Range("D13").Value = "No"
$ws4.Range("A14").Value = "RejectVendorInvoice"
$ws4.Range("B14").Value = "RejectInvoice"
$ws4.Range("C14").Value = "Reject Created Vendor Invoice"
$ws4.Range("D14").Value = "No"

# ---------------------------------------------------------------------
# Selection / active-tab bookkeeping: IND_Regression_FullCycle becomes
# the active sheet/tab (workbookView activeTab 3 -> 2), its selection
# moves to B11, while CHN_Regression_FullCycle's own selection moves to
# D3 and loses tabSelected.
# ---------------------------------------------------------------------
$ws4.Range("D3").Select()
$ws3.Activate()
$ws3.Range("B11").Select()
